$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold, border, centered) from an existing header cell (H1)
# onto the two new header cells so no new style is created.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New column I (I0) and column J (IF) data
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 4

$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 8

$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 7

$ws.Range("I7").Value = 9
$ws.Range("J7").Value = 9

$ws.Range("I8").Value = 7
$ws.Range("J8").Value = 7

$ws.Range("I9").Value = 7
$ws.Range("J9").Value = 7
